$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the static value in C2 with a formula that computes minutes -> hours
$ws.Range("C2").Formula = "=B2/60"
